$wb = $excel.ActiveWorkbook

# Sheet 1: "Bior. Inte. rate of retu."
$ws = $wb.Worksheets.Item(1)
$arr = New-Object 'object[,]' 1,40
$arr[0,0] = 0.1263771079727516
$arr[0,1] = 0.131882303219976
$arr[0,2] = 0.1364141541486866
$arr[0,3] = 0.1405941358863532
$arr[0,4] = 0.1446569484277321
$arr[0,5] = 0.1485416934058388
$arr[0,6] = 0.1524966505124935
$arr[0,7] = 0.1569294417113175
$arr[0,8] = 0.1607816384629603
$arr[0,9] = 0.1643302040077098
$arr[0,10] = 0.1680385073982598
$arr[0,11] = 0.1709874227734602
$arr[0,12] = 0.1739025735896073
$arr[0,13] = 0.1767846392578043
$arr[0,14] = 0.1796337816898061
$arr[0,15] = 0.1824555786624679
$arr[0,16] = 0.1852468843540655
$arr[0,17] = 0.1879817717515009
$arr[0,18] = 0.1907174567173378
$arr[0,19] = 0.1934013512792702
$arr[0,20] = 0.1960853859550391
$arr[0,21] = 0.1987841472119619
$arr[0,22] = 0.2001788398985246
$arr[0,23] = 0.2037646769751738
$arr[0,24] = 0.2054157924359034
$arr[0,25] = 0.2089159433379385
$arr[0,26] = 0.2105441536243976
$arr[0,27] = 0.2140120079238091
$arr[0,28] = 0.215619203768213
$arr[0,29] = 0.2181191135127548
$arr[0,30] = 0.2206170377220371
$arr[0,31] = 0.2231031836858132
$arr[0,32] = 0.2255131861274797
$arr[0,33] = 0.2278718614193154
$arr[0,34] = 0.2302570813832683
$arr[0,35] = 0.2327028854575225
$arr[0,36] = 0.23507120182545
$arr[0,37] = 0.2373724673852868
$arr[0,38] = 0.2392735549595764
$arr[0,39] = 0.2421459389766333
$ws.Range("B2:AO2").Value = $arr

# Sheet 2: "Bior. Biod. prod. cost"
$ws = $wb.Worksheets.Item(2)
$arr = New-Object 'object[,]' 1,40
$arr[0,0] = 4596739.77870624
$arr[0,1] = 6110121.35726426
$arr[0,2] = 7606441.948604696
$arr[0,3] = 9059112.278045349
$arr[0,4] = 10466300.66817262
$arr[0,5] = 11829857.90282767
$arr[0,6] = 13147242.59904813
$arr[0,7] = 14411589.69478941
$arr[0,8] = 15641699.57535029
$arr[0,9] = 16852969.30223284
$arr[0,10] = 18044614.65778005
$arr[0,11] = 19217681.83891526
$arr[0,12] = 20359409.33537341
$arr[0,13] = 21470454.34025065
$arr[0,14] = 22551403.7177843
$arr[0,15] = 23602813.10777451
$arr[0,16] = 24625371.14144303
$arr[0,17] = 25620513.37289041
$arr[0,18] = 26587012.36507841
$arr[0,19] = 27527187.50014611
$arr[0,20] = 28439798.48891104
$arr[0,21] = 29324744.62747439
$arr[0,22] = 30301514.54925343
$arr[0,23] = 31046872.85277091
$arr[0,24] = 31951658.44733849
$arr[0,25] = 32645670.93145756
$arr[0,26] = 33505046.46144887
$arr[0,27] = 34147673.68989405
$arr[0,28] = 34963285.49446175
$arr[0,29] = 35658460.95466016
$arr[0,30] = 36330737.90453882
$arr[0,31] = 36980869.88035478
$arr[0,32] = 37612432.82078218
$arr[0,33] = 38224726.11357758
$arr[0,34] = 38814342.93301651
$arr[0,35] = 39379650.94548032
$arr[0,36] = 39928210.02851827
$arr[0,37] = 40459949.66194697
$arr[0,38] = 40994090.39055227
$arr[0,39] = 41474696.2704834
$ws.Range("B2:AO2").Value = $arr

# Sheet 3: "Bior. Etha. prod. cost"
$ws = $wb.Worksheets.Item(3)
$arr = New-Object 'object[,]' 1,40
$arr[0,0] = 60087728.39640974
$arr[0,1] = 58189816.07497269
$arr[0,2] = 56231476.87163371
$arr[0,3] = 54305140.36340471
$arr[0,4] = 52419394.85982081
$arr[0,5] = 50582863.72008553
$arr[0,6] = 48780563.43513779
$arr[0,7] = 46989557.99154954
$arr[0,8] = 45271286.68733718
$arr[0,9] = 43653081.31838816
$arr[0,10] = 42112638.22048581
$arr[0,11] = 40638156.78463261
$arr[0,12] = 39194255.59569041
$arr[0,13] = 37780355.86970694
$arr[0,14] = 36395794.47106894
$arr[0,15] = 35039922.06510087
$arr[0,16] = 33712327.70067861
$arr[0,17] = 32413535.15625334
$arr[0,18] = 31140655.63685717
$arr[0,19] = 29895340.16626571
$arr[0,20] = 28675080.8877046
$arr[0,21] = 27479002.56999182
$arr[0,22] = 26174456.84941801
$arr[0,23] = 25127877.190469
$arr[0,24] = 23902508.14197984
$arr[0,25] = 22909537.78039474
$arr[0,26] = 21727862.61054406
$arr[0,27] = 20783317.84603306
$arr[0,28] = 19644210.09241179
$arr[0,29] = 18631878.42124644
$arr[0,30] = 17644878.52607512
$arr[0,31] = 16679177.53043213
$arr[0,32] = 15735656.58819508
$arr[0,33] = 14813248.89562073
$arr[0,34] = 13910086.80086249
$arr[0,35] = 13025335.9948332
$arr[0,36] = 12161160.38412445
$arr[0,37] = 11316854.34765258
$arr[0,38] = 10496724.83656216
$arr[0,39] = 9685854.00239254
$ws.Range("B2:AO2").Value = $arr

# Sheet 4: "Bior. Fixed capi. inve."
$ws = $wb.Worksheets.Item(4)
$arr = New-Object 'object[,]' 1,40
$arr[0,0] = 182211002.3322872
$arr[0,1] = 183531423.10008
$arr[0,2] = 184704111.6482563
$arr[0,3] = 185898136.3881662
$arr[0,4] = 187065217.6155345
$arr[0,5] = 188320470.3385141
$arr[0,6] = 189409640.3238878
$arr[0,7] = 189911542.2456904
$arr[0,8] = 190955795.4985475
$arr[0,9] = 191906398.6981884
$arr[0,10] = 192062767.8594966
$arr[0,11] = 192885698.4973143
$arr[0,12] = 193697525.6463881
$arr[0,13] = 194499166.5856293
$arr[0,14] = 195292848.5907481
$arr[0,15] = 196073015.7618885
$arr[0,16] = 196844760.9386395
$arr[0,17] = 197637073.8578999
$arr[0,18] = 198390283.6310282
$arr[0,19] = 199161123.3715741
$arr[0,20] = 199896262.3018367
$arr[0,21] = 200582181.0197221
$arr[0,22] = 201225283.8442827
$arr[0,23] = 201966090.8747699
$arr[0,24] = 202620453.5286523
$arr[0,25] = 203374970.6143662
$arr[0,26] = 204011128.1437108
$arr[0,27] = 204725396.3408416
$arr[0,28] = 205341780.580371
$arr[0,29] = 205989956.4207962
$arr[0,30] = 206628691.171212
$arr[0,31] = 207254477.558488
$arr[0,32] = 207930735.4785731
$arr[0,33] = 208633186.6157412
$arr[0,34] = 209285941.7907197
$arr[0,35] = 209857882.1235547
$arr[0,36] = 210483010.5399592
$arr[0,37] = 211151194.0779106
$arr[0,38] = 212178432.4130722
$arr[0,39] = 212183588.0866981
$ws.Range("B2:AO2").Value = $arr

# Sheet 7: "Biorefinery Steam"
$ws = $wb.Worksheets.Item(7)
$ws.Range("D2").Value = 688822.4568321184
$ws.Range("F2").Value = 674404.9546450275
$ws.Range("J2").Value = 645331.280694496
$ws.Range("K2").Value = 641426.0426595398
$ws.Range("M2").Value = 642382.3327764545
$ws.Range("O2").Value = 643338.2541141368
$ws.Range("Q2").Value = 644293.444535794
$ws.Range("S2").Value = 645249.3114351667
$ws.Range("T2").Value = 645726.7626781734
$ws.Range("X2").Value = 647483.1938594319
$ws.Range("Y2").Value = 648066.432378842
$ws.Range("Z2").Value = 648434.8387509254
$ws.Range("AB2").Value = 649387.4016308556
$ws.Range("AF2").Value = 651287.2117887657
$ws.Range("AH2").Value = 652237.5098450602
$ws.Range("AL2").Value = 654139.2375250295
$ws.Range("AM2").Value = 654615.0421234922

# Sheet 8: "Bior. Cons. elec."
$ws = $wb.Worksheets.Item(8)
$arr = New-Object 'object[,]' 1,40
$arr[0,0] = 48270.4795242657
$arr[0,1] = 48398.09421565657
$arr[0,2] = 48290.04013392414
$arr[0,3] = 48141.63465968824
$arr[0,4] = 47983.5575689352
$arr[0,5] = 47823.14675099129
$arr[0,6] = 47662.24007211855
$arr[0,7] = 47501.00937056764
$arr[0,8] = 47339.71115130379
$arr[0,9] = 47194.28465123341
$arr[0,10] = 47069.56005917418
$arr[0,11] = 46944.79058654243
$arr[0,12] = 46820.32708690343
$arr[0,13] = 46696.90552880164
$arr[0,14] = 46570.39172335657
$arr[0,15] = 46445.70196204122
$arr[0,16] = 46320.66841128656
$arr[0,17] = 46195.5023119673
$arr[0,18] = 46070.72671996498
$arr[0,19] = 45945.57535113927
$arr[0,20] = 45820.69412449987
$arr[0,21] = 45695.76957236964
$arr[0,22] = 45564.03555177029
$arr[0,23] = 45385.58365478595
$arr[0,24] = 45299.29641873537
$arr[0,25] = 45136.0359182326
$arr[0,26] = 45046.33741565012
$arr[0,27] = 44881.64582077137
$arr[0,28] = 44795.05640475154
$arr[0,29] = 44624.45135702083
$arr[0,30] = 44497.13484405542
$arr[0,31] = 44371.23641357909
$arr[0,32] = 44245.39569354206
$arr[0,33] = 44119.32423601466
$arr[0,34] = 43994.17239761062
$arr[0,35] = 43867.86557985593
$arr[0,36] = 43741.62639568256
$arr[0,37] = 43615.14520565924
$arr[0,38] = 43489.76968020611
$arr[0,39] = 43367.5701425345
$ws.Range("B2:AO2").Value = $arr

# Sheet 9: "Biorefinery Excess electricity"
$ws = $wb.Worksheets.Item(9)
$arr = New-Object 'object[,]' 1,40
$arr[0,0] = 208584.8518200926
$arr[0,1] = 216518.0395452966
$arr[0,2] = 225385.6647418447
$arr[0,3] = 234431.8752019143
$arr[0,4] = 243518.5734434969
$arr[0,5] = 252614.4904326223
$arr[0,6] = 261712.444497685
$arr[0,7] = 270811.0576880118
$arr[0,8] = 279909.8118018184
$arr[0,9] = 287659.3307342844
$arr[0,10] = 293651.1047090475
$arr[0,11] = 299642.9414200351
$arr[0,12] = 305634.3337076965
$arr[0,13] = 311625.0764603463
$arr[0,14] = 317618.7767180174
$arr[0,15] = 323610.8068916533
$arr[0,16] = 329602.4588725787
$arr[0,17] = 335595.2525862902
$arr[0,18] = 341587.3375114141
$arr[0,19] = 347579.6005599443
$arr[0,20] = 353571.4514742292
$arr[0,21] = 359563.1374288679
$arr[0,22] = 365624.0125588806
$arr[0,23] = 371627.8369444234
$arr[0,24] = 377624.6263108593
$arr[0,25] = 383609.7103880498
$arr[0,26] = 389613.0834722839
$arr[0,27] = 395599.7351868943
$arr[0,28] = 401599.9194633254
$arr[0,29] = 407640.1765447756
$arr[0,30] = 413635.3508111034
$arr[0,31] = 419629.6078744646
$arr[0,32] = 425623.4476487947
$arr[0,33] = 431617.5359481531
$arr[0,34] = 437610.1247379495
$arr[0,35] = 443605.1685865246
$arr[0,36] = 449599.4460409472
$arr[0,37] = 455593.8289423496
$arr[0,38] = 461571.5504581062
$arr[0,39] = 467114.7482582218
$ws.Range("B2:AO2").Value = $arr
